# daily auto push: 2026-01-21 03:50 UTC
# Insert a new data row for 2026/01/21 03:50 (水, hour=11, rank=201)
# right after the existing 2026/01/21 06:00 entry (row 694), pushing the
# 2026/12/29.. entries (and everything after) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 695..736 down to 696..737, opening up a blank row 695.
$ws.Rows.Item(695).Insert()

# Column A holds a date-like string ("2026/01/21") that must stay a plain
# text value (matching the rest of the column) instead of being
# auto-converted to a real Excel date serial. Force text format first,
# then restore the default "Normal" style so no stray formatting is left
# behind on the cell.
$ws.Range("A695").NumberFormat = "@"
$ws.Range("A695").Value = "2026/01/21"
$ws.Range("A695").Style = "Normal"

$ws.Range("B695").Value = "水"
$ws.Range("C695").Value = 11
$ws.Range("D695").Value = 201
